# Update "想去人数" (want-to-go count) figures in column F across all
# four sheets of the workbook, matching the regenerated gh-pages data
# snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1) # 展览
$ws2 = $wb.Worksheets.Item(2) # 演出
$ws3 = $wb.Worksheets.Item(3) # 本地生活
$ws4 = $wb.Worksheets.Item(4) # 全部类型

# 展览
$ws1.Range("F2").Value = 138
$ws1.Range("F3").Value = 194
$ws1.Range("F4").Value = 453
$ws1.Range("F5").Value = 213
$ws1.Range("F7").Value = 1253
$ws1.Range("F8").Value = 420
$ws1.Range("F10").Value = 60
$ws1.Range("F12").Value = 386
$ws1.Range("F13").Value = 432
$ws1.Range("F14").Value = 806
$ws1.Range("F15").Value = 193
$ws1.Range("F16").Value = 740
$ws1.Range("F17").Value = 300
$ws1.Range("F19").Value = 1039
$ws1.Range("F20").Value = 487
$ws1.Range("F21").Value = 284
$ws1.Range("F22").Value = 98
$ws1.Range("F23").Value = 396
$ws1.Range("F25").Value = 49
$ws1.Range("F26").Value = 488
$ws1.Range("F27").Value = 32

# 演出
$ws2.Range("F4").Value = 372
$ws2.Range("F11").Value = 157
$ws2.Range("F12").Value = 126
$ws2.Range("F13").Value = 35

# 本地生活
$ws3.Range("F2").Value = 351

# 全部类型
$ws4.Range("F2").Value = 351
$ws4.Range("F4").Value = 138
$ws4.Range("F5").Value = 194
$ws4.Range("F6").Value = 453
$ws4.Range("F7").Value = 213
$ws4.Range("F9").Value = 1253
$ws4.Range("F10").Value = 420
$ws4.Range("F13").Value = 60
$ws4.Range("F14").Value = 372
$ws4.Range("F17").Value = 386
$ws4.Range("F20").Value = 432
$ws4.Range("F21").Value = 806
$ws4.Range("F22").Value = 193
$ws4.Range("F23").Value = 740
$ws4.Range("F24").Value = 300
$ws4.Range("F26").Value = 1039
$ws4.Range("F27").Value = 487
$ws4.Range("F30").Value = 284
$ws4.Range("F31").Value = 98
$ws4.Range("F32").Value = 396
$ws4.Range("F34").Value = 157
$ws4.Range("F36").Value = 49
$ws4.Range("F37").Value = 126
$ws4.Range("F38").Value = 35
$ws4.Range("F39").Value = 488
$ws4.Range("F42").Value = 32
